$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.9414315565236273, 0.6976249339774024, 1.3987170663267086, 2.8998042003394526)"
$ws.Range("C2").Value = "JSU(-0.9150914252046711, 1.0702175137455356, 4.351290564771493, 4.947250430425351)"
$ws.Range("D2").Value = "JSU(-0.6829455128984832, 0.9750830952686264, 0.9656932476509564, 2.303694909508288)"
$ws.Range("E2").Value = "JSU(-1.0895870143120048, 1.0830553370980227, 2.90943401634363, 4.8141334387780415)"
